$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level Row1 (E), DP1 (F) source data for rows 2..21 (row index -> (E,F))
$eVals = @(166,209,109,184,178,142,160,145,222,173,222,218,230,174,159,197,180,230,95,160)
$fVals = @(133,84,74,102,121,73,96,103,146,103,144,138,112,110,88,140,128,135,59,104)

# Minute1 (C) is always 20, Second1 (D) is always 0 for every competitor row
$ws.Range("C2:C21").Value = 20
$ws.Range("D2:D21").Value = 0

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
    $ws.Cells.Item($row, 6).Value = $fVals[$i]
}

# Rep1 (G) = Row1 + DP1, entered once over the whole range so Excel stores it
# as a single shared formula (matches the shared si="0" group in the sheet).
$ws.Range("G2").Formula = "=E2+F2"
$ws.Range("G3:G21").Formula = "=E3+F3"

# Stray ";" label added in P16, plus the two new team names that grew the
# roster from 18 to 20 teams.
$ws.Range("P16").Value = ";"
$ws.Range("A20").Value = "LAG 19"
$ws.Range("A21").Value = "LAG 20"

# View state: scroll/zoom and selection moved as part of the same save.
$ws.Range("H17").Select()
$excel.ActiveWindow.Zoom = 115
